$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. GitHub hyperlink display text: leekarensl.github.io -> github.com/leekarensl
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "https://leekarensl.github.io/",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "https://github.com/leekarensl", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Rewrite the "About me" paragraph
#    "My professional background is in client servicing, operations, and
#     project management with key skills in SQL, Excel, Power BI and team
#     management. I am passionate about solving business problems using
#     data so I am constantly learning ... My GitHub portfolio focuses on
#     modelling projects ..."
#    becomes
#    "I have intermediate skills in SQL, Excel, Tableau and Power BI. I am
#     passionate about solving business problems using data, so I am
#     constantly learning ... My GitHub portfolio focuses on data
#     analytics and modelling projects ..."
# ------------------------------------------------------------------

$d.Content.Find.Execute(
    "My professional background is in client servicing, operations, and project management with key skills in SQL, Excel",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "I have intermediate skills in SQL, Excel", 2) | Out-Null

$d.Content.Find.Execute(
    " Power BI and team management. I am passionate about solving business problems using data so I am",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Tableau and Power BI. I am passionate about solving business problems using data, so I am", 2) | Out-Null

$d.Content.Find.Execute(
    "in this area. My GitHub portfolio focuses on modelling",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "in this area. My GitHub portfolio focuses on data analytics and modelling", 2) | Out-Null

# ------------------------------------------------------------------
# 3. Drop the stale rendered-page-break hint that used to sit in front of
#    "Reviewing specifications ..." (Word recomputes these automatically;
#    rewriting the run's text with itself clears the stale hint).
# ------------------------------------------------------------------

$d.Content.Find.Execute(
    "Reviewing specifications for system changes to ensure they match client requirements",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Reviewing specifications for system changes to ensure they match client requirements", 2) | Out-Null

Write-Output "edits applied"
